$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix duplicate company names by appending -N suffixes so the
# "Highlight Duplicate Values" rule below can distinguish repeated
# placements from genuinely duplicated rows. Cells are written in the
# exact order the distinct strings were (re)introduced so the rebuilt
# shared-strings table lines up with the target workbook.
$ws.Range("B3").Value = "AMAZON-2"
$ws.Range("B4").Value = "AMAZON-3"
$ws.Range("B5").Value = "AMAZON-5"
$ws.Range("B11").Value = "ACCOLITE DIGITAL-2"
$ws.Range("B62").Value = "MADHURA GROUP-2"
$ws.Range("B57").Value = "BELCAN-2"
$ws.Range("B54").Value = "TATA CONSULTANCY SERVICES (TCS)-2"
$ws.Range("B52").Value = "INFOSYS-2"
$ws.Range("B46").Value = "ZENSAR-2"
$ws.Range("B43").Value = "VALUELABS-2"
$ws.Range("B42").Value = "COGNIZANT-2"
$ws.Range("B35").Value = "ACCENTURE-2"
$ws.Range("B38").Value = "COGNIZANT-3"
$ws.Range("B41").Value = "CAPGEMINI-2"
$ws.Range("B31").Value = "LARSEN AND TOUBRO INFOTECH (LTI)-2"
$ws.Range("B30").Value = "INFOSYS-3"
$ws.Range("B29").Value = "VIRTUSA-2"
$ws.Range("B23").Value = "INFOR-2"

# Add conditional formatting to highlight duplicate values in column B,
# using Excel's standard "light red fill with dark red text" format.
$rng = $ws.Range("B1:B1048576")
$fc = $rng.FormatConditions.AddUniqueValues()
$fc.DupeUnique = 1
$fc.Font.Color = 393372
$fc.Interior.Color = 13551615

# Move the active selection to B24 (matches the editor's last position).
$null = $ws.Range("B24").Select()
